# Replicates the commit "updated on 5th october":
#   1) The auto-updating "Date Placeholder" field (type datetimeFigureOut) on
#      the slide master and every slide layout is refreshed from 8/25/2022 to
#      9/16/2022 (the value PowerPoint had cached the last time the deck was
#      saved/printed).
#   2) The Title placeholder on slide 4 ("Xpath") is nudged down and made a
#      little shorter (its autofit box tightened around the now-2-line title).

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes, $newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
            # ppPlaceholderDate == 16: the "datetimeFigureOut" field shape.
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "9/16/2022"

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes $newDate

# Every layout under the master has its own copy of the placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes $newDate
}

# Slide 4 ("Xpath") title box: move down/shrink slightly (x & width unchanged).
$slide4 = $p.Slides.Item(4)
$title = $slide4.Shapes.Item(1)
$title.Top = 33.86771653543307
$title.Height = 48.738977477952766
